$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 9.736532333333333
$ws.Cells.Item(2, 8).Value = 29.209597
$ws.Cells.Item(2, 9).Value = 0.3545698647072128
$ws.Cells.Item(2, 10).Value = 0.3545698647072129
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.1433513333333333
$ws.Cells.Item(2, 14).Value = 0.430054
$ws.Cells.Item(2, 15).Value = 0.002710896760593916
$ws.Cells.Item(2, 16).Value = 0.002710896760593916
$ws.Cells.Item(2, 17).Value = 1.395744892026444
$ws.Cells.Item(2, 18).Value = 12.561704028238
$ws.Cells.Item(2, 19).Value = 0.0009612022976390062
$ws.Cells.Item(2, 20).Value = 0.0009612022976390064
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 9.736532333333333
$ws.Cells.Item(3, 8).Value = 29.209597
$ws.Cells.Item(3, 9).Value = 0.3545698647072128
$ws.Cells.Item(3, 10).Value = 0.3545698647072129
$ws.Cells.Item(3, 13).Value = 43.12631833333334
$ws.Cells.Item(3, 15).Value = 0.8155556976531461
$ws.Cells.Item(3, 16).Value = 0.8155556976531461
$ws.Cells.Item(3, 17).Value = 419.9007928701261
$ws.Cells.Item(3, 18).Value = 3779.107135831136
$ws.Cells.Item(3, 19).Value = 0.2891714733780725
$ws.Cells.Item(3, 20).Value = 0.2891714733780726
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 9.736532333333333
$ws.Cells.Item(4, 8).Value = 29.209597
$ws.Cells.Item(4, 9).Value = 0.3545698647072128
$ws.Cells.Item(4, 10).Value = 0.3545698647072129
$ws.Cells.Item(4, 13).Value = 9.610002999999999
$ws.Cells.Item(4, 14).Value = 28.830009
$ws.Cells.Item(4, 15).Value = 0.18173340558626
$ws.Cells.Item(4, 16).Value = 0.1817334055862599
$ws.Cells.Item(4, 17).Value = 93.56810493293032
$ws.Cells.Item(4, 18).Value = 842.1129443963729
$ws.Cells.Item(4, 19).Value = 0.06443718903150122
$ws.Cells.Item(4, 20).Value = 0.06443718903150122
$ws.Cells.Item(5, 7).Value = 10.17625966666667
$ws.Cells.Item(5, 9).Value = 0.37058316962423
$ws.Cells.Item(5, 10).Value = 0.37058316962423
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.1433513333333333
$ws.Cells.Item(5, 14).Value = 0.430054
$ws.Cells.Item(5, 15).Value = 0.002710896760593916
$ws.Cells.Item(5, 16).Value = 0.002710896760593916
$ws.Cells.Item(5, 17).Value = 1.458780391562889
$ws.Cells.Item(5, 18).Value = 13.129023524066
$ws.Cells.Item(5, 19).Value = 0.001004612714064951
$ws.Cells.Item(5, 20).Value = 0.001004612714064951
$ws.Cells.Item(6, 7).Value = 10.17625966666667
$ws.Cells.Item(6, 9).Value = 0.37058316962423
$ws.Cells.Item(6, 10).Value = 0.37058316962423
$ws.Cells.Item(6, 13).Value = 43.12631833333334
$ws.Cells.Item(6, 15).Value = 0.8155556976531461
$ws.Cells.Item(6, 16).Value = 0.8155556976531461
$ws.Cells.Item(6, 17).Value = 438.8646138273273
$ws.Cells.Item(6, 18).Value = 3949.781524445946
$ws.Cells.Item(6, 19).Value = 0.3022312154414031
$ws.Cells.Item(6, 20).Value = 0.3022312154414031
$ws.Cells.Item(7, 7).Value = 10.17625966666667
$ws.Cells.Item(7, 9).Value = 0.37058316962423
$ws.Cells.Item(7, 10).Value = 0.37058316962423
$ws.Cells.Item(7, 13).Value = 9.610002999999999
$ws.Cells.Item(7, 14).Value = 28.830009
$ws.Cells.Item(7, 15).Value = 0.18173340558626
$ws.Cells.Item(7, 16).Value = 0.1817334055862599
$ws.Cells.Item(7, 17).Value = 97.79388592544565
$ws.Cells.Item(7, 18).Value = 880.1449733290109
$ws.Cells.Item(7, 19).Value = 0.06734734146876196
$ws.Cells.Item(7, 20).Value = 0.06734734146876196
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 7.547331666666667
$ws.Cells.Item(8, 8).Value = 22.641995
$ws.Cells.Item(8, 9).Value = 0.2748469656685572
$ws.Cells.Item(8, 10).Value = 0.2748469656685572
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 0.1433513333333333
$ws.Cells.Item(8, 14).Value = 0.430054
$ws.Cells.Item(8, 15).Value = 0.002710896760593916
$ws.Cells.Item(8, 16).Value = 0.002710896760593916
$ws.Cells.Item(8, 17).Value = 1.081920057525556
$ws.Cells.Item(8, 18).Value = 9.737280517730001
$ws.Cells.Item(8, 19).Value = 0.0007450817488899588
$ws.Cells.Item(8, 20).Value = 0.0007450817488899588
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 7.547331666666667
$ws.Cells.Item(9, 8).Value = 22.641995
$ws.Cells.Item(9, 9).Value = 0.2748469656685572
$ws.Cells.Item(9, 10).Value = 0.2748469656685572
$ws.Cells.Item(9, 13).Value = 43.12631833333334
$ws.Cells.Item(9, 15).Value = 0.8155556976531461
$ws.Cells.Item(9, 16).Value = 0.8155556976531461
$ws.Cells.Item(9, 17).Value = 325.488628023914
$ws.Cells.Item(9, 18).Value = 2929.397652215226
$ws.Cells.Item(9, 19).Value = 0.2241530088336704
$ws.Cells.Item(9, 20).Value = 0.2241530088336704
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 7.547331666666667
$ws.Cells.Item(10, 8).Value = 22.641995
$ws.Cells.Item(10, 9).Value = 0.2748469656685572
$ws.Cells.Item(10, 10).Value = 0.2748469656685572
$ws.Cells.Item(10, 13).Value = 9.610002999999999
$ws.Cells.Item(10, 14).Value = 28.830009
$ws.Cells.Item(10, 15).Value = 0.18173340558626
$ws.Cells.Item(10, 16).Value = 0.1817334055862599
$ws.Cells.Item(10, 17).Value = 72.52987995866167
$ws.Cells.Item(10, 18).Value = 652.7689196279549
$ws.Cells.Item(10, 19).Value = 0.04994887508599677
$ws.Cells.Item(10, 20).Value = 0.04994887508599676
